# Update internal test data on the "none" sheet: pivot the per-date rows
# (A3:E5 holding ref/measure_name/comment/date/value) into a single wide
# row with one date column per period (matching the layout already used
# on the "week"/"month" sheets).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("none")

# Drop the now-redundant rows 3-5 (their date/value pairs move into new
# columns F/G on row 2 below).
$ws.Range("A3:A5").EntireRow.Delete() | Out-Null

# Row 1: D1/E1 used to hold the text labels "date"/"value"; they - and two
# new columns F1/G1 - now hold the actual period-end dates.
$ws.Range("D1").Value = 43890
$ws.Range("E1").Value = 44511
$ws.Range("F1").Value = 44662
$ws.Range("G1").Value = 44834
$ws.Range("D1:G1").NumberFormat = "m/d/yy"

# Row 2: the four measurement values (previously spread across D2:D5/E2:E5)
# now sit side by side as whole numbers.
$ws.Range("D2").Value = 60
$ws.Range("E2").Value = 620
$ws.Range("F2").Value = 151
$ws.Range("G2").Value = 172
$ws.Range("D2:G2").NumberFormat = "0"

# Widen the new date columns to fit their contents.
$ws.Columns.Item(6).ColumnWidth = 12.833333333333334
$ws.Columns.Item(7).ColumnWidth = 14.5
$ws.Columns.Item(8).ColumnWidth = 15.166666666666666

# Leave the selection on the first data row, as in the saved file.
$ws.Range("A2").Select() | Out-Null

# Nudge the saved window position down slightly (best effort - matches the
# author's last on-screen Excel window position).
$wb.Windows.Item(1).Top = 1800
